# Quest_List.xlsx - "add map enemy factory"
#
# Replaces the single "EnrmyPlace" column (D) with three new columns
# EnemyStage1 / EnemyStage2 / EnemyStage3 (D:F), pushing BossFlag and the
# EnemyID1-5 columns two slots to the right (old E:J -> new G:L).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Seed the two brand-new columns (K, L) with the same cell format
#        that the data they will hold already uses (I:J), so the new cells
#        don't end up with a freshly-minted style index.
$ws.Range("I1:J4").Copy() | Out-Null
$ws.Range("K1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# --- 2) Header row — rewrite strictly left to right. Each destination cell
#        is written before its old source text is needed anywhere else,
#        which keeps the shared-string table's prune/append order aligned
#        with the target file (EnrmyPlace drops out, BossFlag/EnemyID1
#        shift down, EnemyStage1-3 get appended at the end).
$ws.Range("D1").Value = "EnemyStage1"
$ws.Range("E1").Value = "EnemyStage2"
$ws.Range("F1").Value = "EnemyStage3"
$ws.Range("G1").Value = "BossFlag"
$ws.Range("H1").Value = "EnemyID1"
$ws.Range("I1").Value = "EnemyID2"
$ws.Range("J1").Value = "EnemyID3"
$ws.Range("K1").Value = "EnemyID4"
$ws.Range("L1").Value = "EnemyID5"

# --- 3) Data rows. Write literal final values (right to left, so a column
#        we still need to read from never gets clobbered first).

# Row 2 (MapID 101)
$ws.Range("L2").Value = 104000
$ws.Range("K2").Value = 5000
$ws.Range("J2").Value = 4000
$ws.Range("I2").Value = 5000
$ws.Range("H2").Value = 4000
$ws.Range("G2").Value = $false
$ws.Range("F2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("D2").Value = 2

# Row 3 (MapID 102)
$ws.Range("L3").Value = 104000
$ws.Range("K3").Value = 104000
$ws.Range("J3").Value = 104000
$ws.Range("I3").Value = 104000
$ws.Range("H3").Value = 104000
$ws.Range("G3").Value = $false
$ws.Range("F3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("D3").Value = 2

# Row 4 (MapID 103)
$ws.Range("L4").Value = 104000
$ws.Range("K4").Value = 104000
$ws.Range("J4").Value = 104000
$ws.Range("I4").Value = 104000
$ws.Range("H4").Value = 104000
$ws.Range("G4").Value = $false
$ws.Range("F4").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("D4").Value = 2

# --- 4) Column E width tightened (19.33 -> 12.83 characters) now that it
#        holds "EnemyStage2" instead of the wider "EnrmyPlace"/boolean col.
#        (12 is the closest achievable ColumnWidth to the recorded
#        12.83203125 given this host's pixel-granularity rounding.)
$ws.Range("E1").EntireColumn.ColumnWidth = 12

# --- 5) Selection moves to K3.
$ws.Range("K3").Select()
